# Actualización automática 2025-08-19 14:15:09
#
# Insert a new client "WONG SANCHEZ PAULA SOFIA" (advisor "GUERRERO FAREZ
# FABIAN MAURICIO") as a new row just above "ZUÑIGA CORONEL MARCIA LUZMILA"
# on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing the
# existing totals/summary row down by one, and bumping the "de 55" -> "de 56"
# counters on the "VENTAS POR GRUPO" summary row (the advisor now has 56
# clients instead of 55).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R, data rows 2-56, summary 57)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new blank row at 56 (this shifts the old row 56 "ZUÑIGA CORONEL
# MARCIA LUZMILA" down to row 57, and the old summary row 57 down to row 58).
$ws1.Rows.Item(56).Insert()

# Fill in the new row 56 with the new client.
$ws1.Cells.Item(56, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws1.Cells.Item(56, 2).Value = "WONG SANCHEZ PAULA SOFIA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(56, $col).Value = 0
}

# Update the "X de 55" counters (now row 58) to "X de 56".
$null = $ws1.Range("C58:R58").Replace("de 55", "de 56")

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G, data rows 2-56, totals 57)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(56).Insert()

$ws2.Cells.Item(56, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws2.Cells.Item(56, 2).Value = "WONG SANCHEZ PAULA SOFIA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(56, $col).Value = 0
}
